$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark "Absent" (column H) = 1 for every attendance date row (3-18)
$ws.Range("H3:H18").Value = 1

# For the dates where the student was also counted in Total Attendance / Real
# (rows 10, 13 and 17), set columns D and E to 1 as well
$ws.Range("D10:E10").Value = 1
$ws.Range("D13:E13").Value = 1
$ws.Range("D17:E17").Value = 1
